$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39 (ALC) -- diff @@ -2564
$ws.Range("H39").Value = 119.166664
$ws.Range("I39").Value = 65.333336
$ws.Range("J39").Value = 173
$ws.Range("K39").Value = 196.000008
$ws.Range("L39").Value = 519
$ws.Range("M39").Value = 99.99999199999999
$ws.Range("N39").Value = -1111

# Row 69 (ALC) -- diff @@ -4052
$ws.Range("H69").Value = 15000
$ws.Range("J69").Value = 15000
$ws.Range("L69").Value = 45000
$ws.Range("N69").Value = -46748

# Row 70 (ALC) -- diff @@ -4098
$ws.Range("H70").Value = 11300.2
$ws.Range("I70").Value = 500
$ws.Range("K70").Value = 1500
$ws.Range("M70").Value = -1230

# Row 72 (ALC) -- diff @@ -4193
$ws.Range("H72").Value = 15000
$ws.Range("J72").Value = 15000
$ws.Range("L72").Value = 135000
$ws.Range("N72").Value = -143736

# Row 73 (ALC) -- diff @@ -4239
$ws.Range("H73").Value = 11300.2
$ws.Range("I73").Value = 500
$ws.Range("K73").Value = 1500
$ws.Range("M73").Value = -564

# Row 114 (ALC) -- diff @@ -6302
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = $null
$ws.Range("N114").Value = 0

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM) -- diff @@ -7786
$ws.Range("H2").Value = 4816.423
$ws.Range("I2").Value = 1304.7646
$ws.Range("K2").Value = 1304.7646
$ws.Range("M2").Value = -1191.7646

# Row 32 (ARM) -- diff @@ -9277
$ws.Range("H32").Value = 1525936.2
$ws.Range("I32").Value = 2904.0308
$ws.Range("J32").Value = 13900573
$ws.Range("K32").Value = 2904.0308
$ws.Range("L32").Value = 13900573
$ws.Range("M32").Value = -2617.0308
$ws.Range("N32").Value = -13901147

# Row 116 (ARM) -- diff @@ -13318
$ws.Range("H116").Value = 4816.423
$ws.Range("I116").Value = 1304.7646
$ws.Range("K116").Value = 1304.7646
$ws.Range("M116").Value = 989.2354

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM) -- diff @@ -14714
$ws.Range("H3").Value = 4816.423
$ws.Range("I3").Value = 1304.7646
$ws.Range("K3").Value = 1304.7646
$ws.Range("M3").Value = -1190.7646

# Row 107 (BSM) -- diff @@ -19780
$ws.Range("H107").Value = 3341290.8
$ws.Range("I107").Value = 5269248
$ws.Range("J107").Value = 11182.546
$ws.Range("K107").Value = 5269248
$ws.Range("L107").Value = 11182.546
$ws.Range("M107").Value = -5267328
$ws.Range("N107").Value = -15022.546

# Row 134 (BSM) -- diff @@ -21049
$ws.Range("H134").Value = 1287223.8
$ws.Range("I134").Value = 1616996.8
$ws.Range("K134").Value = 4850990.4
$ws.Range("M134").Value = -4848455.4

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (CRP) -- diff @@ -22533
$ws.Range("H22").Value = 1624915.5
$ws.Range("J22").Value = 3510
$ws.Range("L22").Value = 3510
$ws.Range("N22").Value = -4210

# Row 31 (CRP) -- diff @@ -22980
$ws.Range("H31").Value = 4735
$ws.Range("I31").Value = 1245.75
$ws.Range("K31").Value = 1245.75
$ws.Range("M31").Value = -950.75

# Row 34 (CRP) -- diff @@ -23133
$ws.Range("H34").Value = 4735
$ws.Range("I34").Value = 1245.75
$ws.Range("K34").Value = 1245.75
$ws.Range("M34").Value = -1043.75

# Row 58 (CRP) -- diff @@ -24291
$ws.Range("H58").Value = 30315254
$ws.Range("I58").Value = 41673892
$ws.Range("J58").Value = 25552.777
$ws.Range("K58").Value = 41673892
$ws.Range("L58").Value = 25552.777
$ws.Range("M58").Value = -41673689
$ws.Range("N58").Value = -25958.777

# Row 86 (CRP) -- diff @@ -25642
$ws.Range("H86").Value = 13708.444
$ws.Range("J86").Value = 19092.25
$ws.Range("L86").Value = 19092.25
$ws.Range("N86").Value = -21338.25

# Row 89 (CRP) -- diff @@ -25789
$ws.Range("H89").Value = 13708.444
$ws.Range("J89").Value = 19092.25
$ws.Range("L89").Value = 95461.25
$ws.Range("N89").Value = -106693.25

# Row 136 (CRP) -- diff @@ -28050
$ws.Range("H136").Value = 30315254
$ws.Range("I136").Value = 41673892
$ws.Range("J136").Value = 25552.777
$ws.Range("K136").Value = 125021676
$ws.Range("L136").Value = 76658.33099999999
$ws.Range("M136").Value = -125019126
$ws.Range("N136").Value = -81758.33099999999

$ws = $wb.Worksheets.Item("CUL")
# Row 116 (CUL) -- diff @@ -34195
$ws.Range("H116").Value = 6713.9287
$ws.Range("I116").Value = 6917.0835
$ws.Range("J116").Value = 5495
$ws.Range("K116").Value = 20751.2505
$ws.Range("L116").Value = 16485
$ws.Range("M116").Value = -17309.2505
$ws.Range("N116").Value = -23369

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (GSM) -- diff @@ -40954
$ws.Range("H113").Value = 8825.6
$ws.Range("I113").Value = 4549
$ws.Range("K113").Value = 4549
$ws.Range("M113").Value = -2379

# Row 126 (GSM) -- diff @@ -41573
$ws.Range("H126").Value = 62515264
$ws.Range("I126").Value = 100004424
$ws.Range("K126").Value = 300013272
$ws.Range("M126").Value = -300010802

# Row 132 (GSM) -- diff @@ -41864
$ws.Range("H132").Value = 29414964
$ws.Range("I132").Value = 33336376
$ws.Range("J132").Value = 4368.75
$ws.Range("K132").Value = 100009128
$ws.Range("L132").Value = 13106.25
$ws.Range("M132").Value = -100006598
$ws.Range("N132").Value = -18166.25

$ws = $wb.Worksheets.Item("LTW")
# Row 43 (LTW) -- diff @@ -44457
$ws.Range("H43").Value = 2407001.2
$ws.Range("I43").Value = 35006
$ws.Range("J43").Value = 3000000
$ws.Range("K43").Value = 35006
$ws.Range("L43").Value = 3000000
$ws.Range("M43").Value = -34813
$ws.Range("N43").Value = -3000386

# Row 61 (LTW) -- diff @@ -45348
$ws.Range("H61").Value = 8362.684999999999
$ws.Range("I61").Value = 7813.2
$ws.Range("J61").Value = 8973.223
$ws.Range("K61").Value = 7813.2
$ws.Range("L61").Value = 8973.223
$ws.Range("M61").Value = -7611.2
$ws.Range("N61").Value = -9377.223

# Row 113 (LTW) -- diff @@ -47863
$ws.Range("H113").Value = 8362.684999999999
$ws.Range("I113").Value = 7813.2
$ws.Range("J113").Value = 8973.223
$ws.Range("K113").Value = 7813.2
$ws.Range("L113").Value = 8973.223
$ws.Range("M113").Value = -5643.2
$ws.Range("N113").Value = -13313.223

# Row 122 (LTW) -- diff @@ -48289
$ws.Range("H122").Value = 4250.5
$ws.Range("I122").Value = 2751.25
$ws.Range("J122").Value = 6249.5
$ws.Range("K122").Value = 8253.75
$ws.Range("L122").Value = 18748.5
$ws.Range("M122").Value = -5803.75
$ws.Range("N122").Value = -23648.5

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (WVR) -- diff @@ -52264
$ws.Range("H62").Value = 11299.2
$ws.Range("I62").Value = 17624
$ws.Range("K62").Value = 17624
$ws.Range("M62").Value = -17000

# Row 65 (WVR) -- diff @@ -52408
$ws.Range("H65").Value = 11299.2
$ws.Range("I65").Value = 17624
$ws.Range("K65").Value = 88120
$ws.Range("M65").Value = -85000

# Row 107 (WVR) -- diff @@ -54445
$ws.Range("H107").Value = 2967.2354
$ws.Range("I107").Value = 2904.3076
$ws.Range("J107").Value = 3171.75
$ws.Range("K107").Value = 8712.9228
$ws.Range("L107").Value = 9515.25
$ws.Range("M107").Value = -6792.9228
$ws.Range("N107").Value = -13355.25

# Row 126 (WVR) -- diff @@ -55355
$ws.Range("H126").Value = 2706.48
$ws.Range("J126").Value = 6207.5713
$ws.Range("L126").Value = 18622.7139
$ws.Range("N126").Value = -23562.7139

# Row 132 (WVR) -- diff @@ -55646
$ws.Range("H132").Value = 7968.0234
$ws.Range("I132").Value = 5065.448
$ws.Range("K132").Value = 15196.344
$ws.Range("M132").Value = -12666.344
